# PeerReviewRS1.docx edit:
#  - Q1 paragraph ("1.  What is the main argument ...") gains a trailing
#    sentence: " Underline where you find it."
#  - Q3 paragraph ("3. What are the main points ...") gains a trailing
#    sentence: " Underline where these points are stated."
#  - The "_GoBack" bookmark, which previously lived alone in its own empty
#    paragraph right after Q1, is removed from there and re-created at the
#    very end of the Q3 paragraph (after the newly appended sentence).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Append " Underline where you find it." to the Q1 paragraph.
# ---------------------------------------------------------------------
$q1 = $d.Content
$q1.Find.Execute("What is the main argument of the article being summarized?", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$q1.Collapse(0)                       # wdCollapseEnd
$q1Pos = $q1.Start

$q1Insert = $d.Range($q1Pos, $q1Pos)
$q1Insert.InsertAfter(" Underline where you find it.")

# Match the surrounding "Garamond" body font on the freshly inserted run.
$d.Range($q1Pos, $q1Insert.End).Font.Name = "Garamond"

# ---------------------------------------------------------------------
# 2. Remove the "_GoBack" bookmark from its current (otherwise-empty)
#    paragraph right after Q1.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 3. Append " Underline where these points are stated." to the Q3
#    paragraph, then re-add the "_GoBack" bookmark collapsed at the new
#    end of that paragraph (right after the just-inserted sentence).
# ---------------------------------------------------------------------
$q3 = $d.Content
$q3.Find.Execute("What are the main points of the article being summarized?", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$q3.Collapse(0)                       # wdCollapseEnd
$q3Pos = $q3.Start

$q3Insert = $d.Range($q3Pos, $q3Pos)
$q3Insert.InsertAfter(" Underline where these points are stated.")
$q3End = $q3Insert.End

# Match the surrounding "Garamond" body font on the freshly inserted run.
$d.Range($q3Pos, $q3End).Font.Name = "Garamond"

# Adding a bookmark collapsed exactly at a paragraph-end position is
# unreliable, so nudge it into place: insert a throw-away character after
# the target spot (so the spot is no longer paragraph-final), drop the
# bookmark there, then delete the throw-away character again. The
# bookmark sticks to the boundary and ends up collapsed at the true end.
$placeholder = $d.Range($q3End, $q3End)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($q3End, $q3End)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($q3End, $q3End + 1).Delete()

Write-Output "Done."
